# Reformat national phone numbers in column E by inserting a space after
# the first 3 digits of the national number (area code), e.g.
#   "+57 3164071898"  -> "+57 316 4071898"
#   "+1 7543045130"   -> "+1 754 3045130"
# and fill in the missing "Departamento" (state/province) for the two
# Cucuta-based tenants (rows 30-31) with "Santander".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$phoneUpdates = @{
    "E2"  = "+57 316 4071898"
    "E3"  = "+57 316 6841377"
    "E4"  = "+1 754 3045130"
    "E5"  = "+1 804 5546462"
    "E6"  = "+1 754 3045130"
    "E7"  = "+57 315 6074411"
    "E9"  = "+57 314 6327654"
    "E11" = "+57 318 6490366"
    "E12" = "+57 316 8703511"
    "E13" = "+57 310 6431401"
    "E14" = "+57 317 8002340"
    "E15" = "+57 315 4630661"
    "E18" = "+57 317 4899555"
    "E19" = "+57 316 8712962"
    "E20" = "+57 316 6418549"
    "E21" = "+57 317 4286020"
    "E22" = "+57 316 5551566"
    "E23" = "+57 311 8209865"
    "E24" = "+57 310 40882215"
    "E25" = "+57 301 2541700"
    "E28" = "+57 316 4186372"
    "E29" = "*57 312 8696103"
    "E30" = "+57 320 8961514"
    "E31" = "+57 315 6389025"
    "E32" = "+57 300 6087604"
    "E33" = "+57 300 6087520"
}

foreach ($ref in $phoneUpdates.Keys) {
    $ws.Range($ref).Value = $phoneUpdates[$ref]
}

# Fill in the state/province ("Departamento") for the Cucuta rows that
# were previously missing it.
$ws.Range("H30").Value = "Santander"
$ws.Range("H31").Value = "Santander"

# Match the author's last on-screen selection from the commit.
$ws.Range("H6:I6").Select()
